$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 29 == "Google Shape;48;p1" — the problem-statement question textbox.
$sh = $s.Shapes.Item(29)

# --- Resize / reposition the textbox (it grew taller to fit the new wording) ---
$sh.Top = 42.590625
$sh.Height = 74.35244140625

# --- Replace the body text with the revised question ---
$tr = $sh.TextFrame.TextRange
$tr.Text = ""
$main = $tr.InsertAfter("Can a neural network or multiclassification model be trained to identify crops planted in a field with at least 90% accuracy, so that Prime Agri" + [char]8217 + "s management team can make a decision on the most competitive seed line to pursue in the new territory within the next year? ")
$main.LanguageID = "en-US"

$mainLen = $main.Length - 1
$boldPart = $main.Characters(1, $mainLen)
$boldPart.Font.Bold = $true
